$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.17532746519896136
$ws.Range("B1").Value = 0.17478542411376452
$ws.Range("A2").Value = -0.090409751280960293
$ws.Range("B2").Value = 0.089055210751265435
$ws.Range("A3").Value = 0.013885254835461325
$ws.Range("B3").Value = -0.014333616831073925
$ws.Range("A4").Value = -0.1956663819990041
$ws.Range("B4").Value = 0.19432983705323537
$ws.Range("A5").Value = -0.18832983748679588
$ws.Range("B5").Value = 0.18561520809508281
$ws.Range("A6").Value = -0.084680217291768667
$ws.Range("B6").Value = 0.084583332312029746
$ws.Range("A7").Value = -0.064583332843980656
$ws.Range("B7").Value = 0.064376021294940244
$ws.Range("A8").Value = -0.016236839274888837
$ws.Range("B8").Value = 0.016215962802120565
$ws.Range("A9").Value = -0.010215963261382299
$ws.Range("B9").Value = 0.010203801016814396
$ws.Range("A10").Value = -0.0042038014771890175
$ws.Range("B10").Value = 0.0042046052709068249
$ws.Range("A11").Value = 0.00029539427706737342
$ws.Range("B11").Value = -0.0002977943006641226
$ws.Range("A12").Value = 0.0062977938402046796
$ws.Range("B12").Value = -0.0063303314191403004
$ws.Range("A13").Value = 0.012330330959453129
$ws.Range("B13").Value = -0.0123508583852594
$ws.Range("A14").Value = 0.024350857892411426
$ws.Range("B14").Value = -0.024405892032828724
$ws.Range("A15").Value = -0.02105167182947465
$ws.Range("B15").Value = 0.021026993404406724
$ws.Range("A16").Value = -0.015026993863574312
$ws.Range("B16").Value = 0.015004324153656157
$ws.Range("A17").Value = -0.0090043246147333278
$ws.Range("B17").Value = 0.0089999995215999107
$ws.Range("A18").Value = -0.074375673277319976
$ws.Range("B18").Value = 0.074285296897674158
$ws.Range("A19").Value = -0.027096920379364509
$ws.Range("B19").Value = 0.02701371667370589
$ws.Range("A20").Value = -0.018013717117465688
$ws.Range("B20").Value = 0.018004305484389249
$ws.Range("A21").Value = -0.0090043059287117089
$ws.Range("B21").Value = 0.0089999995552707546
$ws.Range("A22").Value = -0.093940466030350223
$ws.Range("B22").Value = 0.093629808500903522
$ws.Range("A23").Value = -0.08462980894831329
$ws.Range("B23").Value = 0.084125876212950423
$ws.Range("A24").Value = -0.042125876853677191
$ws.Range("B24").Value = 0.041999999355874174
$ws.Range("A25").Value = -0.094934991254760348
$ws.Range("B25").Value = 0.094690075668211193
$ws.Range("A26").Value = -0.088690076121608286
$ws.Range("B26").Value = 0.088374217722403614
$ws.Range("A27").Value = -0.082374218178177916
$ws.Range("B27").Value = 0.081292939728166314
$ws.Range("A28").Value = -0.075292940194056079
$ws.Range("B28").Value = 0.074543736287114193
$ws.Range("A29").Value = -0.062543736793609028
$ws.Range("B29").Value = 0.062172318557200512
$ws.Range("A30").Value = -0.042172319112065981
$ws.Range("B30").Value = 0.042020047335702415
$ws.Range("A31").Value = -0.027020047867274855
$ws.Range("B31").Value = 0.027000806859939885
$ws.Range("A32").Value = -0.0060008074261448385
$ws.Range("B32").Value = 0.0059999995180826105

$ws.Columns.Item(1).ColumnWidth = 14.8
$ws.Columns.Item(2).ColumnWidth = 15.6
